{"js": "// Replace the math-problem answers in the table with the new values.\n// Each old value is unique in the document, so a simple search/replace\n// per pair is unambiguous.\nconst replacements = [\n  [\"19\u00d716=304\", \"87\u00d712=1044\"],\n  [\"54\u00d777=4158\", \"98\u00d748=4704\"],\n  [\"40\u00d777=3080\", \"47\u00d779=3713\"],\n  [\"49\u00d725=1225\", \"14\u00d742=588\"],\n  [\"30\u00d729=870\", \"67\u00d793=6231\"],\n  [\"44\u00d711=484\", \"19\u00d721=399\"],\n  [\"42\u00d783=3486\", \"90\u00d773=6570\"],\n  [\"18\u00d752=936\", \"21\u00d788=1848\"],\n  [\"56\u00d749=2744\", \"43\u00d776=3268\"],\n  [\"16\u00d719=304\", \"42\u00d741=1722\"],\n  [\"87\u00d743=3741\", \"68\u00d782=5576\"],\n  [\"32\u00d739=1248\", \"21\u00d768=1428\"],\n  [\"76\u00d720=1520\", \"83\u00d721=1743\"],\n  [\"60\u00d769=4140\", \"23\u00d749=1127\"],\n  [\"41\u00d745=1845\", \"99\u00d769=6831\"],\n  [\"98\u00d768=6664\", \"98\u00d734=3332\"],\n  [\"22\u00d788=1936\", \"60\u00d722=1320\"],\n  [\"55\u00d716=880\", \"59\u00d769=4071\"],\n  [\"46\u00d769=3174\", \"97\u00d734=3298\"],\n  [\"49\u00d795=4655\", \"71\u00d713=923\"],\n  [\"63\u00d719=1197\", \"14\u00d761=854\"],\n  [\"59\u00d720=1180\", \"39\u00d746=1794\"],\n  [\"84\u00d786=7224\", \"33\u00d744=1452\"],\n  [\"89\u00d799=8811\", \"42\u00d792=3864\"],\n  [\"87\u00d722=1914\", \"71\u00d722=1562\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the math-problem answers in the table with the new values.\n# Each old value is unique in the document, so a simple Find/Replace\n# per pair is unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"19\u00d716=304\", \"87\u00d712=1044\"),\n    @(\"54\u00d777=4158\", \"98\u00d748=4704\"),\n    @(\"40\u00d777=3080\", \"47\u00d779=3713\"),\n    @(\"49\u00d725=1225\", \"14\u00d742=588\"),\n    @(\"30\u00d729=870\", \"67\u00d793=6231\"),\n    @(\"44\u00d711=484\", \"19\u00d721=399\"),\n    @(\"42\u00d783=3486\", \"90\u00d773=6570\"),\n    @(\"18\u00d752=936\", \"21\u00d788=1848\"),\n    @(\"56\u00d749=2744\", \"43\u00d776=3268\"),\n    @(\"16\u00d719=304\", \"42\u00d741=1722\"),\n    @(\"87\u00d743=3741\", \"68\u00d782=5576\"),\n    @(\"32\u00d739=1248\", \"21\u00d768=1428\"),\n    @(\"76\u00d720=1520\", \"83\u00d721=1743\"),\n    @(\"60\u00d769=4140\", \"23\u00d749=1127\"),\n    @(\"41\u00d745=1845\", \"99\u00d769=6831\"),\n    @(\"98\u00d768=6664\", \"98\u00d734=3332\"),\n    @(\"22\u00d788=1936\", \"60\u00d722=1320\"),\n    @(\"55\u00d716=880\", \"59\u00d769=4071\"),\n    @(\"46\u00d769=3174\", \"97\u00d734=3298\"),\n    @(\"49\u00d795=4655\", \"71\u00d713=923\"),\n    @(\"63\u00d719=1197\", \"14\u00d761=854\"),\n    @(\"59\u00d720=1180\", \"39\u00d746=1794\"),\n    @(\"84\u00d786=7224\", \"33\u00d744=1452\"),\n    @(\"89\u00d799=8811\", \"42\u00d792=3864\"),\n    @(\"87\u00d722=1914\", \"71\u00d722=1562\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards, $false, $false, $find.Forward, $find.Wrap, $false, $find.Replacement.Text, 2)  # wdReplaceAll\n}\n"}
